$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-format style (s="2") from A343 down to the new date cells A344:A357
$ws.Range("A343").Copy()
$ws.Range("A344:A357").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(344,1).Value = 44418
$ws.Cells.Item(344,2).Value = 11
$ws.Cells.Item(344,3).Value = 43
$ws.Cells.Item(344,4).Value = 168.3040432110846

$ws.Cells.Item(345,1).Value = 44419
$ws.Cells.Item(345,2).Value = 3
$ws.Cells.Item(345,3).Value = 44
$ws.Cells.Item(345,4).Value = 172.2180907276214

$ws.Cells.Item(346,1).Value = 44420
$ws.Cells.Item(346,2).Value = 14
$ws.Cells.Item(346,3).Value = 51
$ws.Cells.Item(346,4).Value = 199.6164233433794

$ws.Cells.Item(347,1).Value = 44421
$ws.Cells.Item(347,2).Value = 5
$ws.Cells.Item(347,3).Value = 45
$ws.Cells.Item(347,4).Value = 176.1321382441583

$ws.Cells.Item(348,1).Value = 44422
$ws.Cells.Item(348,2).Value = 7
$ws.Cells.Item(348,3).Value = 50
$ws.Cells.Item(348,4).Value = 195.7023758268425

$ws.Cells.Item(349,1).Value = 44423
$ws.Cells.Item(349,2).Value = 2
$ws.Cells.Item(349,3).Value = 48
$ws.Cells.Item(349,4).Value = 187.8742807937688

$ws.Cells.Item(350,1).Value = 44424
$ws.Cells.Item(350,2).Value = 2
$ws.Cells.Item(350,3).Value = 44
$ws.Cells.Item(350,4).Value = 172.2180907276214

$ws.Cells.Item(351,1).Value = 44425
$ws.Cells.Item(351,2).Value = 4
$ws.Cells.Item(351,3).Value = 37
$ws.Cells.Item(351,4).Value = 144.8197581118635

$ws.Cells.Item(352,1).Value = 44426
$ws.Cells.Item(352,2).Value = 1
$ws.Cells.Item(352,3).Value = 35
$ws.Cells.Item(352,4).Value = 136.9916630787898

$ws.Cells.Item(353,1).Value = 44427
$ws.Cells.Item(353,2).Value = 2
$ws.Cells.Item(353,3).Value = 23
$ws.Cells.Item(353,4).Value = 90.02309288034756

$ws.Cells.Item(354,1).Value = 44428
$ws.Cells.Item(354,2).Value = 2
$ws.Cells.Item(354,3).Value = 20
$ws.Cells.Item(354,4).Value = 78.28095033073701

$ws.Cells.Item(355,1).Value = 44429
$ws.Cells.Item(355,2).Value = 4
$ws.Cells.Item(355,3).Value = 17
$ws.Cells.Item(355,4).Value = 66.53880778112647

$ws.Cells.Item(356,1).Value = 44430
$ws.Cells.Item(356,2).Value = 3
$ws.Cells.Item(356,3).Value = 18
$ws.Cells.Item(356,4).Value = 70.45285529766332

$ws.Cells.Item(357,1).Value = 44431
$ws.Cells.Item(357,2).Value = 4
$ws.Cells.Item(357,3).Value = 20
$ws.Cells.Item(357,4).Value = 78.28095033073701
